# "add vega popup for CH" — Travel History workbook edit
#
# Semantic changes (decoded from the OOXML diff):
#   1. Row 20, column C (Subdivisions) for the Switzerland/Geneva entry:
#        "Geneva" -> "Genève"
#   2. Row 25, column C (Subdivisions) for the Switzerland/Lucerne entry:
#        "Lucerne" -> "Luzern"
#   3. A brand-new row 45 is appended: Switzerland / CHE / Fribourg / 2025-02
#   4. The active selection moves from E19 to E21 (and the viewport scrolls
#      so row 16 is the top-left visible row — best-effort via Select()).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the Geneva subdivision to its native-language spelling.
$ws.Range("C20").Value = "Genève"

# 2) Rename the Lucerne subdivision to its native-language spelling.
$ws.Range("C25").Value = "Luzern"

# 3) Append a new travel-history row for Fribourg (Switzerland).
$ws.Range("A45").Value = "Switzerland"
$ws.Range("B45").Value = "CHE"
$ws.Range("C45").Value = "Fribourg"
$ws.Range("D45").Value = "2025-02"

# 4) Update the saved selection to match the author's final cursor position.
$ws.Range("E21").Select()
